# ------------------------------------------------------------------
# Refresh the Market Board snapshot values (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job leve tables.
# Mirrors the scheduled-runner data refresh described in the commit.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1499.5
$ws.Range("I62").Value = 1499.5
$ws.Range("K62").Value = 1499.5
$ws.Range("M62").Value = -875.5
$ws.Range("H65").Value = 1499.5
$ws.Range("I65").Value = 1499.5
$ws.Range("K65").Value = 7497.5
$ws.Range("M65").Value = -4377.5
$ws.Range("H107").Value = 4886.9585
$ws.Range("I107").Value = 3225.7368
$ws.Range("K107").Value = 3225.7368
$ws.Range("M107").Value = -1305.7368
$ws.Range("H116").Value = 52649.8
$ws.Range("J116").Value = 54499.75
$ws.Range("L116").Value = 54499.75
$ws.Range("N116").Value = -61383.75
$ws.Range("H132").Value = 2795.6667
$ws.Range("I132").Value = 2868.4546
$ws.Range("K132").Value = 8605.363799999999
$ws.Range("M132").Value = -6075.363799999999
$ws.Range("H137").Value = 6899793
$ws.Range("I137").Value = 10528315
$ws.Range("K137").Value = 31584945
$ws.Range("M137").Value = -31582395

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14839594
$ws.Range("I61").Value = 6063086
$ws.Range("K61").Value = 6063086
$ws.Range("M61").Value = -6062874
$ws.Range("H74").Value = 1557563.2
$ws.Range("I74").Value = 2721168.8
$ws.Range("K74").Value = 2721168.8
$ws.Range("M74").Value = -2720294.8
$ws.Range("H77").Value = 1557563.2
$ws.Range("I77").Value = 2721168.8
$ws.Range("K77").Value = 13605844
$ws.Range("M77").Value = -13601476
$ws.Range("H122").Value = 2821.5
$ws.Range("I122").Value = 2611.8
$ws.Range("K122").Value = 7835.400000000001
$ws.Range("M122").Value = -5385.400000000001
$ws.Range("H132").Value = 5943.9375
$ws.Range("I132").Value = 4370.4736
$ws.Range("J132").Value = 6974.8276
$ws.Range("K132").Value = 13111.4208
$ws.Range("L132").Value = 20924.4828
$ws.Range("M132").Value = -10581.4208
$ws.Range("N132").Value = -25984.4828
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 14839594
$ws.Range("I136").Value = 6063086
$ws.Range("K136").Value = 18189258
$ws.Range("M136").Value = -18186708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 375.625
$ws.Range("I80").Value = 540.5
$ws.Range("K80").Value = 540.5
$ws.Range("M80").Value = 457.5
$ws.Range("H83").Value = 375.625
$ws.Range("I83").Value = 540.5
$ws.Range("K83").Value = 2702.5
$ws.Range("M83").Value = 2289.5
$ws.Range("H134").Value = 6947904.5
$ws.Range("I134").Value = 6947058.5
$ws.Range("K134").Value = 20841175.5
$ws.Range("M134").Value = -20838640.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 752388.5600000001
$ws.Range("I31").Value = 1116463.5
$ws.Range("J31").Value = 4012.1667
$ws.Range("K31").Value = 1116463.5
$ws.Range("L31").Value = 4012.1667
$ws.Range("M31").Value = -1116168.5
$ws.Range("N31").Value = -4602.1667
$ws.Range("H34").Value = 752388.5600000001
$ws.Range("I34").Value = 1116463.5
$ws.Range("J34").Value = 4012.1667
$ws.Range("K34").Value = 1116463.5
$ws.Range("L34").Value = 4012.1667
$ws.Range("M34").Value = -1116261.5
$ws.Range("N34").Value = -4416.1667
$ws.Range("H58").Value = 11676912
$ws.Range("I58").Value = 16667875
$ws.Range("K58").Value = 16667875
$ws.Range("M58").Value = -16667672
$ws.Range("H99").Value = 12677.305
$ws.Range("J99").Value = 7144.222
$ws.Range("L99").Value = 7144.222
$ws.Range("N99").Value = -10140.222
$ws.Range("H126").Value = 12677.305
$ws.Range("J126").Value = 7144.222
$ws.Range("L126").Value = 21432.666
$ws.Range("N126").Value = -26372.666
$ws.Range("H132").Value = 4197
$ws.Range("I132").Value = 4398.375
$ws.Range("J132").Value = 3874.8
$ws.Range("K132").Value = 13195.125
$ws.Range("L132").Value = 11624.4
$ws.Range("M132").Value = -10665.125
$ws.Range("N132").Value = -16684.4
$ws.Range("H134").Value = 5169.972
$ws.Range("I134").Value = 3458.111
$ws.Range("J134").Value = 5740.593
$ws.Range("K134").Value = 10374.333
$ws.Range("L134").Value = 17221.779
$ws.Range("M134").Value = -7839.332999999999
$ws.Range("N134").Value = -22291.779
$ws.Range("H136").Value = 11676912
$ws.Range("I136").Value = 16667875
$ws.Range("K136").Value = 50003625
$ws.Range("M136").Value = -50001075

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 212.83333
$ws.Range("I2").Value = 65.5
$ws.Range("J2").Value = 286.5
$ws.Range("K2").Value = 393
$ws.Range("L2").Value = 1719
$ws.Range("M2").Value = -280
$ws.Range("N2").Value = -1945
$ws.Range("H5").Value = 2359364.5
$ws.Range("I5").Value = 1701136.5
$ws.Range("K5").Value = 5103409.5
$ws.Range("M5").Value = -5103297.5
$ws.Range("H38").Value = 161.42105
$ws.Range("J38").Value = 73.22221999999999
$ws.Range("L38").Value = 219.66666
$ws.Range("N38").Value = -913.66666
$ws.Range("H68").Value = 3989.3835
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 4257.179
$ws.Range("K68").Value = 2997
$ws.Range("L68").Value = 12771.537
$ws.Range("M68").Value = -2186
$ws.Range("N68").Value = -14393.537
$ws.Range("H71").Value = 3989.3835
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 4257.179
$ws.Range("K71").Value = 8991
$ws.Range("L71").Value = 38314.611
$ws.Range("M71").Value = -4935
$ws.Range("N71").Value = -46426.611
$ws.Range("H86").Value = 333.3889
$ws.Range("I86").Value = 429.6
$ws.Range("K86").Value = 1288.8
$ws.Range("M86").Value = -102.8000000000002
$ws.Range("H89").Value = 333.3889
$ws.Range("I89").Value = 429.6
$ws.Range("K89").Value = 3866.4
$ws.Range("M89").Value = 2061.6
$ws.Range("H107").Value = 3620.35
$ws.Range("I107").Value = 694.2857
$ws.Range("J107").Value = 5195.923
$ws.Range("K107").Value = 2082.8571
$ws.Range("L107").Value = 15587.769
$ws.Range("M107").Value = -162.8571000000002
$ws.Range("N107").Value = -19427.769
$ws.Range("H113").Value = 666
$ws.Range("J113").Value = 755.44446
$ws.Range("L113").Value = 2266.33338
$ws.Range("N113").Value = -6606.33338
$ws.Range("H127").Value = 5990.294
$ws.Range("J127").Value = 5990.294
$ws.Range("L127").Value = 17970.882
$ws.Range("N127").Value = -27890.882
$ws.Range("H135").Value = 2359364.5
$ws.Range("I135").Value = 1701136.5
$ws.Range("K135").Value = 15310228.5
$ws.Range("M135").Value = -15307693.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12938.174
$ws.Range("I70").Value = 14953
$ws.Range("K70").Value = 14953
$ws.Range("M70").Value = -14683
$ws.Range("H73").Value = 12938.174
$ws.Range("I73").Value = 14953
$ws.Range("K73").Value = 14953
$ws.Range("M73").Value = -14017
$ws.Range("H126").Value = 4492
$ws.Range("I126").Value = 3984.5
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 11953.5
$ws.Range("L126").Value = 14998.5
$ws.Range("H127").Value = 92198.836
$ws.Range("J127").Value = 92198.836
$ws.Range("L127").Value = 92198.836
$ws.Range("N127").Value = -102118.836
$ws.Range("H132").Value = 16258.37
$ws.Range("I132").Value = 11387.105
$ws.Range("K132").Value = 34161.315
$ws.Range("M132").Value = -31631.315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2023.75
$ws.Range("I16").Value = 1031.6666
$ws.Range("K16").Value = 1031.6666
$ws.Range("M16").Value = -861.6666
$ws.Range("H40").Value = 6169.3335
$ws.Range("I40").Value = 5504.25
$ws.Range("J40").Value = 7499.5
$ws.Range("K40").Value = 5504.25
$ws.Range("L40").Value = 7499.5
$ws.Range("M40").Value = -5368.25
$ws.Range("N40").Value = -7771.5
$ws.Range("H136").Value = 15974173
$ws.Range("I136").Value = 15627001
$ws.Range("J136").Value = 16668516
$ws.Range("K136").Value = 46881003
$ws.Range("L136").Value = 50005548
$ws.Range("M136").Value = -46878453
$ws.Range("N136").Value = -50010648

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 105738.3
$ws.Range("I81").Value = 5571.5713
$ws.Range("K81").Value = 11143.1426
$ws.Range("M81").Value = -10082.1426
$ws.Range("H84").Value = 105738.3
$ws.Range("I84").Value = 5571.5713
$ws.Range("K84").Value = 55715.713
$ws.Range("M84").Value = -50411.713
$ws.Range("H122").Value = 40082.125
$ws.Range("I122").Value = 3718.963
$ws.Range("K122").Value = 11156.889
$ws.Range("M122").Value = -8706.889000000001
$ws.Range("H132").Value = 11115441
$ws.Range("I132").Value = 16670702
$ws.Range("K132").Value = 50012106
$ws.Range("M132").Value = -50009576
$ws.Range("H141").Value = 72357.5
$ws.Range("J141").Value = 72357.5
$ws.Range("L141").Value = 72357.5
$ws.Range("N141").Value = -82717.5
